$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new value to A2 (new line under existing data)
$ws.Range("A2").Value = 6789

# Move/select cell A4, matching the selection recorded in the saved file
$ws.Range("A4").Select()
